$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.146.17"
$ws.Range("E2").Value = "  +0.24%  "

$ws.Range("D3").Value = "1.903.97"
$ws.Range("E3").Value = "  +0.73%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'305.96"
$ws.Range("E5").Value = "  -0.51%  "

$ws.Range("D7").Value = "'0.5224"
$ws.Range("E7").Value = "  +1.59%  "

$ws.Range("D8").Value = "'0.3755"
$ws.Range("E8").Value = "  +0.40%  "

$ws.Range("D9").Value = "'0.07249"
$ws.Range("E9").Value = "  +0.51%  "

$ws.Range("E10").Value = "  -0.12%  "

$ws.Range("D11").Value = "'0.9026"

$ws.Range("D12").Value = "'0.08536"
$ws.Range("E12").Value = "  +11.73%  "

$ws.Range("D13").Value = "1.904.92"
$ws.Range("E13").Value = "  +0.77%  "

$ws.Range("D14").Value = "'95.10"
$ws.Range("E14").Value = "  +0.41%  "

$ws.Range("D15").Value = "'5.289"
$ws.Range("E15").Value = "  +0.33%  "

$ws.Range("E16").Value = "  -0.01%  "

$ws.Range("D17").Value = "'0.000008631"
$ws.Range("E17").Value = "  +1.82%  "

$ws.Range("E18").Value = "  +0.58%  "

$ws.Range("E19").Value = "  +0.03%  "

$ws.Range("D20").Value = "27.182.80"
$ws.Range("E20").Value = "  +0.25%  "

$ws.Range("E21").Value = "  -0.09%  "

$ws.Range("D22").Value = "2.146.97"
$ws.Range("E22").Value = "  +1.02%  "

$ws.Range("E23").Value = "  +0.48%  "

$ws.Range("D24").Value = "'6.424"
$ws.Range("E24").Value = "  +0.20%  "

$ws.Range("E25").Value = "  +3.35%  "

$ws.Range("D26").Value = "'146.94"
$ws.Range("E26").Value = "  +0.61%  "

$ws.Range("E27").Value = "  -1.97%  "

$ws.Range("E28").Value = "  +0.68%  "

$ws.Range("D29").Value = "'115.00"
$ws.Range("E29").Value = "  +0.35%  "

$ws.Range("D30").Value = "'4.813"
$ws.Range("E30").Value = "  -0.80%  "

$ws.Range("E31").Value = "  -1.13%  "

$ws.Range("E32").Value = "  +0.75%  "

$ws.Range("D33").Value = "'0.8063"
$ws.Range("E33").Value = "  +4.95%  "

$ws.Range("D34").Value = "'0.05045"
$ws.Range("E34").Value = "  -0.79%  "

$ws.Range("E35").Value = "  -0.09%  "

$ws.Range("E36").Value = "  +4.94%  "

$ws.Range("D37").Value = "'2.943"
$ws.Range("E37").Value = "  -1.08%  "

$ws.Range("D38").Value = "'2.617"
$ws.Range("E38").Value = "  +0.52%  "

$ws.Range("D39").Value = "'0.5719"
$ws.Range("E39").Value = "  +2.10%  "

$ws.Range("E40").Value = "  -0.19%  "

$ws.Range("E41").Value = "  -0.07%  "

$ws.Range("D42").Value = "'9.080"
$ws.Range("E42").Value = "  +1.26%  "

$ws.Range("D43").Value = "'6.631"
$ws.Range("E43").Value = "  -0.26%  "

$ws.Range("D44").Value = "'115.90"
$ws.Range("E44").Value = "  -1.78%  "

$ws.Range("D45").Value = "'0.1517"
$ws.Range("E45").Value = "  +0.16%  "

$ws.Range("D46").Value = "'0.4863"
$ws.Range("E46").Value = "  +1.27%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'10.16"
$ws.Range("E47").Value = "  -0.50%  "

$ws.Range("B48").Value = "PaxDollar"
$ws.Range("C48").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D48").Value = "'0.9999"
$ws.Range("E48").Value = "  +0.05%  "

$ws.Range("D49").Value = "'1.614"
$ws.Range("E49").Value = "  +1.51%  "

$ws.Range("D50").Value = "'37.50"
$ws.Range("E50").Value = "  +0.15%  "

$ws.Range("D51").Value = "'64.01"
$ws.Range("E51").Value = "  +0.13%  "
